# Update values in result_data_KNN.xlsx (commit: "Update Name of Algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A3"  = -21.932
    "E3"  = 16.325
    "A21" = -20.24
    "A23" = -20.53599999999999
    "E24" = 16.841
    "A25" = -21.664
    "B27" = 5.835
    "B31" = 6.075
    "B39" = 7.594000000000001
    "B48" = 5.274
    "B51" = 6.236
    "B52" = 5.399
    "A53" = -21.993
    "B55" = 5.044
    "B56" = 4.976999999999999
    "A57" = -22.053
    "B57" = 5.415999999999999
    "E57" = 16.486
    "A59" = -22.5
    "E61" = 16.468
    "A69" = -21.649
    "E70" = 17.708
    "B73" = 6.803999999999999
    "A79" = -21.192
    "A83" = -22.006
    "E86" = 16.345
    "B89" = 5.678
    "B90" = 5.767
    "A93" = -21.476
    "E98" = 16.368
    "E100" = 16.609
    "E102" = 16.446
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
